# Apply updated "want to go" (F column) counts across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 538
$ws1.Range("F4").Value  = 190
$ws1.Range("F6").Value  = 498
$ws1.Range("F7").Value  = 99
$ws1.Range("F9").Value  = 41
$ws1.Range("F10").Value = 6648
$ws1.Range("F11").Value = 228
$ws1.Range("F12").Value = 362
$ws1.Range("F13").Value = 2875
$ws1.Range("F14").Value = 185
$ws1.Range("F15").Value = 319
$ws1.Range("F16").Value = 258
$ws1.Range("F17").Value = 528

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 11

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 11
$ws4.Range("F5").Value  = 538
$ws4.Range("F6").Value  = 190
$ws4.Range("F8").Value  = 498
$ws4.Range("F9").Value  = 99
$ws4.Range("F11").Value = 41
$ws4.Range("F13").Value = 6648
$ws4.Range("F15").Value = 228
$ws4.Range("F16").Value = 362
$ws4.Range("F17").Value = 2875
$ws4.Range("F18").Value = 185
$ws4.Range("F19").Value = 319
$ws4.Range("F20").Value = 258
$ws4.Range("F21").Value = 528
